$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 19 with the benchmarking data for the fixedTimeStep / symmetry_coarse case
$ws.Range("B19").Value = "shippingBox_convection_symmetry_coarse_fixedTimeStep"
$ws.Range("D19").Value = "fixed(0,1)"
$ws.Range("E19").Value = 45775
$ws.Range("F19").Value = "1351s"
$ws.Range("G19").Value = 271.12112999999999
$ws.Range("G19").NumberFormat = "#,##0.000000"
$ws.Range("H19").Value = 274.73237999999998
$ws.Range("J19").Value = "Time = 6000s"

# Update the selection to match the saved workbook state
$ws.Range("J20").Select()
